$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Comments" header cell picks up the sheet's alternate (but visually
# identical) Normal formatting once the column starts being used below.
$ws.Range("M1").NumberFormat = "General"

# Row 14 (n = 10) - replace existing trial values, clear the "comment" note,
# and add an Average formula in column L.
$ws.Range("B14").Value = 0.0000367165
$ws.Range("C14").Value = 0.0000293255
$ws.Range("D14").Value = 0.0000240803
$ws.Range("E14").Value = 0.0000236034
$ws.Range("F14").Value = 0.0000197887
$ws.Range("G14").Value = 0.0000181198
$ws.Range("H14").Value = 0.0000157356
$ws.Range("I14").Value = 0.0000190735
$ws.Range("J14").Value = 0.0000171661
$ws.Range("K14").Value = 0.0000171661
$ws.Range("L14").Formula = "=AVERAGE(B14:K14)"
$ws.Range("M14").ClearContents()
$ws.Range("M14").NumberFormat = "General"

# Row 15 (n = 100)
$ws.Range("B15").Value = 0.0002672672
$ws.Range("C15").Value = 0.0002603531
$ws.Range("D15").Value = 0.0002534389
$ws.Range("E15").Value = 0.0002527237
$ws.Range("F15").Value = 0.0002527237
$ws.Range("G15").Value = 0.0002577305
$ws.Range("H15").Value = 0.0002548695
$ws.Range("I15").Value = 0.0002617836
$ws.Range("J15").Value = 0.0002522469
$ws.Range("K15").Value = 0.0002596378
$ws.Range("L15").Formula = "=AVERAGE(B15:K15)"

# Row 16 (n = 1000)
$ws.Range("B16").Value = 0.0189437866
$ws.Range("C16").Value = 0.0185317993
$ws.Range("D16").Value = 0.0185585022
$ws.Range("E16").Value = 0.0186450481
$ws.Range("F16").Value = 0.0185697079
$ws.Range("G16").Value = 0.0185995102
$ws.Range("H16").Value = 0.0186738968
$ws.Range("I16").Value = 0.0186161995
$ws.Range("J16").Value = 0.0185492039
$ws.Range("K16").Value = 0.0200004578
$ws.Range("L16").Formula = "=AVERAGE(B16:K16)"

# Row 17 (n = 10000)
$ws.Range("B17").Value = 1.90036726
$ws.Range("C17").Value = 1.9290881157
$ws.Range("D17").Value = 1.9752721786
$ws.Range("E17").Value = 2.0134222507
$ws.Range("F17").Value = 1.9555208683
$ws.Range("G17").Value = 1.9984002113
$ws.Range("H17").Value = 1.9526746273
$ws.Range("I17").Value = 1.9439589977
$ws.Range("J17").Value = 1.9800419807
$ws.Range("K17").Value = 1.9200298786
$ws.Range("L17").Formula = "=AVERAGE(B17:K17)"

# Row 18 (n = 100000)
$ws.Range("B18").Value = 520.0134701729
$ws.Range("C18").Value = 463.6133487225
$ws.Range("D18").Value = 381.1040716171
$ws.Range("E18").Value = 567.3169505596
$ws.Range("F18").Value = 486.6429803371
$ws.Range("G18").Value = 509.2805206776
$ws.Range("H18").Value = 543.4538860321
$ws.Range("I18").Value = 496.919598341
$ws.Range("J18").Value = 525.3025047779
$ws.Range("K18").Value = 515.0943915844
$ws.Range("L18").Formula = "=AVERAGE(B18:K18)"

# Update the active selection to match the final cursor position.
$ws.Range("C21").Select() | Out-Null
